$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray (label-less) row 13 that only held the professor's
# id/name value; deleting it shifts every following row up by one so the
# column-A labels land on the rows the target layout expects
# (dimension A1:C23 instead of A1:C24).
$ws.Rows("13:13").Delete()

# After the shift, several data cells (columns B/C) need their text swapped
# out for different content so the sheet matches the target state; the
# column-A labels already line up correctly after the row delete.
$ws.Range("B10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C10").Value = "6495737 - Durval Rodrigues Junior"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B18").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C18").Value = "6495737 - Durval Rodrigues Junior"

$ws.Range("B19").Value = "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo."
$ws.Range("C19").Value = "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo."

$ws.Range("B20").Value = "Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3"
$ws.Range("C20").Value = "Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3"

$ws.Range("B21").Value = "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# "01/01/2012" looks like a date, so a plain .Value assignment would get
# auto-converted to a date serial number (and pick up a brand-new
# number-format style). Route it through a text formula and paste-special
# "values only" instead, which keeps it a literal shared-string cell using
# the existing column style, exactly like the rest of the sheet.
$ws.Range("B15").Formula = "=""01/01/2012"""
$ws.Range("C15").Formula = "=""01/01/2012"""
$ws.Range("B15:C15").Copy()
$ws.Range("B15:C15").PasteSpecial(-4163)
